# Insert a new data row before the current row 90. This pushes the old
# row 90..165 down to 91..166 (Excel automatically shifts all cell
# contents and formatting down, extending the sheet from A1:R165 to
# A1:R166), matching a new weekly price observation being recorded for
# "Vega Modelo de Temuco" / "Bruselas (repollito)".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("90:90").Insert()

$ws.Cells.Item(90, 1).Value  = 10
$ws.Cells.Item(90, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(90, 3).Value  = "La Araucanía"
$ws.Cells.Item(90, 4).Value  = 45062
$ws.Cells.Item(90, 5).Value  = 9
$ws.Cells.Item(90, 6).Value  = 100112035
$ws.Cells.Item(90, 7).Value  = "Bruselas (repollito)"
$ws.Cells.Item(90, 8).Value  = "Sin especificar"
$ws.Cells.Item(90, 9).Value  = "Primera"
$ws.Cells.Item(90, 10).Value = 35
$ws.Cells.Item(90, 11).Value = 28000
$ws.Cells.Item(90, 12).Value = 28000
$ws.Cells.Item(90, 13).Value = 28000
$ws.Cells.Item(90, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(90, 15).Value = "Región Metropolitana"
$ws.Cells.Item(90, 16).Value = 1867
$ws.Cells.Item(90, 17).Value = 15
$ws.Cells.Item(90, 18).Value = "Hortaliza"
